# Append a new data row (row 85) to the "Data" sheet, mirroring the
# layout/formatting of the existing rows (date/time in column A, integer
# counters in columns B:O).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row   # xlUp
$newRow = $lastRow + 1

# Copy the whole previous row (values + formatting) down one row, then
# overwrite the copied values with the new data set.
$ws.Range("A" + $lastRow + ":O" + $lastRow).Copy()
$ws.Range("A" + $newRow + ":O" + $newRow).PasteSpecial(-4104)  # xlPasteAll
$excel.CutCopyMode = $false

$ws.Cells.Item($newRow, 1).Value = 45720.93949074074
$ws.Cells.Item($newRow, 2).Value = 10
$ws.Cells.Item($newRow, 3).Value = 6
$ws.Cells.Item($newRow, 4).Value = 235
$ws.Cells.Item($newRow, 5).Value = 411
$ws.Cells.Item($newRow, 6).Value = 387
$ws.Cells.Item($newRow, 7).Value = 474
$ws.Cells.Item($newRow, 8).Value = 3381
$ws.Cells.Item($newRow, 9).Value = 474
$ws.Cells.Item($newRow, 10).Value = 2026
$ws.Cells.Item($newRow, 11).Value = 208
$ws.Cells.Item($newRow, 12).Value = 415
$ws.Cells.Item($newRow, 13).Value = 30
$ws.Cells.Item($newRow, 14).Value = 3612
$ws.Cells.Item($newRow, 15).Value = 4604
